$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Level 4 alien row, mirrors the structure of rows 2-4.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "`"SkeletalMesh'/Game/Assets/3D/Yoppoppo_Generic/Mesh_Generic.Mesh_Generic'`""
$ws.Range("C5").Value = "`"AnimBlueprint'/Game/Assets/3D/Yoppoppo_Generic/AnimationBP_Generic.AnimationBP_Generic'`""
$ws.Range("D5").Value = "`"DataTable'/Game/Assets/AlienDatatables/Level1-4_AlienData.Level1-4_AlienData'`""
$ws.Range("E5").Value = "A Yoppoppo traveler arrives, looking somewhat pained. It seems to be in a rush!`n`"Skroog... Weddi Skroog?`""

$ws.Range("E5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 60

$ws.Range("E5").Select() | Out-Null
